# Apply the "update excel and db schemas" edit:
# In the "Abstract pages" sheet, the "Group link" row beneath each of the
# three event-group pages ("All events", "All summer events",
# "All winter events") is expanded from a single cross-link into a
# three-line list that also references the other two sibling group pages.
# This pushes the three following tables ("All points of interest",
# "All itineraries", "All services") down by six rows in total, and makes
# "Abstract pages" the active sheet/tab.

$wb = $excel.ActiveWorkbook
$wsTables = $wb.Worksheets.Item("Content tables")
$ws = $wb.Worksheets.Item("Abstract pages")

# ---------------------------------------------------------------------
# 1) Insert new rows (bottom-most insertion point first, so the row
#    numbers used below for the still-untouched parts of the sheet stay
#    valid while we work).
# ---------------------------------------------------------------------

# --- "All winter events" page (row 149 = "Group link" row) ---
$ws.Range("A150:B151").EntireRow.Insert()
$ws.Range("B149").Value = "All events"
$ws.Range("B150").Value = "Summer events"
$ws.Range("B151").Value = "Page for the kind of topic event  (thumbnail + name + date + short description)"

# --- "All summer events" page (row 136 = "Group link" row) ---
$ws.Range("A137:B138").EntireRow.Insert()
$ws.Range("B136").Value = "All events"
$ws.Range("B137").Value = "Winter events"
$ws.Range("B138").Value = "Page for the kind of topic event  (thumbnail + name + date + short description)"

# --- "All events" page (row 123 = "Group link" row) ---
$ws.Range("A124:B125").EntireRow.Insert()
$ws.Range("B123").Value = "Summer events"
$ws.Range("B124").Value = "Winter events"
$ws.Range("B125").Value = "Page for the kind of topic event (thumbnail + name + date + short description)"

# ---------------------------------------------------------------------
# 2) Resize / reposition the affected Excel Tables to match their new,
#    final extents (row inserts do not auto-resize tables here).
# ---------------------------------------------------------------------

$ws.ListObjects.Item("Table1930").Resize($ws.Range("A110:B125"))   # All events
$ws.ListObjects.Item("Table2031").Resize($ws.Range("A127:B140"))   # All summer events
$ws.ListObjects.Item("Table2132").Resize($ws.Range("A142:B155"))   # All winter events
$ws.ListObjects.Item("Table2233").Resize($ws.Range("A157:B168"))   # All points of interest
$ws.ListObjects.Item("Table2334").Resize($ws.Range("A170:B181"))   # All itineraries
$ws.ListObjects.Item("Table2435").Resize($ws.Range("A183:B194"))   # All services

# ---------------------------------------------------------------------
# 3) Update sheet view / active tab so "Abstract pages" becomes the
#    selected, visible sheet, scrolled/selected the same way as the
#    authored edit.
# ---------------------------------------------------------------------

$wsTables.Select()
$wsTables.Range("A3").Select()

$ws.Select()
$ws.Application.ActiveWindow.ScrollRow = 129
$ws.Range("A154").Select()
